$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force the Price/Volume columns to Text format while writing the
# new values so Excel does not auto-convert numeric-looking strings (e.g. "1.00",
# "125.00", "54.357.54") into numbers and silently drop trailing zeros / dot
# separators. Afterwards restore General format + Normal style so the cells
# keep no explicit style, matching their original (unstyled) state.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = '54.357.54'
$ws.Range("E2").Value = '  -6.30%  '
$ws.Range("D3").Value = '2.878.32'
$ws.Range("E3").Value = '  -9.34%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '468.64'
$ws.Range("E5").Value = '  -11.79%  '
$ws.Range("D6").Value = '125.00'
$ws.Range("E6").Value = '  -6.86%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '2.876.23'
$ws.Range("E8").Value = '  -9.34%  '
$ws.Range("E9").Value = '  -11.34%  '
$ws.Range("E10").Value = '  -8.73%  '
$ws.Range("D11").Value = '0.0954'
$ws.Range("E11").Value = '  -14.46%  '
$ws.Range("E12").Value = '  -17.36%  '
$ws.Range("E13").Value = '  -4.90%  '
$ws.Range("D14").Value = '3.367.89'
$ws.Range("E14").Value = '  -9.53%  '
$ws.Range("D15").Value = '23.04'
$ws.Range("E15").Value = '  -10.04%  '
$ws.Range("D16").Value = '54.326.24'
$ws.Range("E16").Value = '  -6.56%  '
$ws.Range("D17").Value = '2.873.96'
$ws.Range("E17").Value = '  -9.68%  '
$ws.Range("E18").Value = '  -14.53%  '
$ws.Range("D19").Value = '5.27'
$ws.Range("E19").Value = '  -9.57%  '
$ws.Range("D20").Value = '11.29'
$ws.Range("E20").Value = '  -14.41%  '
$ws.Range("D21").Value = '7.01'
$ws.Range("E21").Value = '  -12.98%  '
$ws.Range("D22").Value = '295.06'
$ws.Range("E22").Value = '  -17.57%  '
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("D24").Value = '0.435'
$ws.Range("E24").Value = '  -15.61%  '
$ws.Range("D25").Value = '58.22'
$ws.Range("E25").Value = '  -16.36%  '
$ws.Range("E26").Value = '  +0.44%  '
$ws.Range("E27").Value = '  -10.78%  '
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("D29").Value = '0.0₃0795'
$ws.Range("E29").Value = '  -16.10%  '
$ws.Range("D30").Value = '6.05'
$ws.Range("E30").Value = '  -12.28%  '
$ws.Range("D31").Value = '1.10'
$ws.Range("E31").Value = '  -7.79%  '
$ws.Range("D32").Value = '6.08'
$ws.Range("E32").Value = '  -12.47%  '
$ws.Range("E33").Value = '  -15.52%  '
$ws.Range("D34").Value = '18.50'
$ws.Range("E34").Value = '  -14.43%  '
$ws.Range("D35").Value = '138.03'
$ws.Range("E35").Value = '  -13.95%  '
$ws.Range("D36").Value = '4.11'
$ws.Range("E36").Value = '  -16.51%  '
$ws.Range("D37").Value = '5.34'
$ws.Range("E37").Value = '  -14.74%  '
$ws.Range("E38").Value = '  -15.75%  '
$ws.Range("D39").Value = '22.84'
$ws.Range("E39").Value = '  -12.25%  '
$ws.Range("D40").Value = '2.902.99'
$ws.Range("E40").Value = '  -9.41%  '
$ws.Range("D41").Value = '0.0609'
$ws.Range("E41").Value = '  -13.13%  '
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("D43").Value = '35.07'
$ws.Range("E43").Value = '  -13.38%  '
$ws.Range("E44").Value = '  -14.70%  '
$ws.Range("D45").Value = '0.934'
$ws.Range("E45").Value = '  -13.89%  '
$ws.Range("E46").Value = '  -11.94%  '
$ws.Range("D47").Value = '3.35'
$ws.Range("E47").Value = '  -16.14%  '
$ws.Range("D48").Value = '2.030.19'
$ws.Range("E48").Value = '  -10.72%  '
$ws.Range("D49").Value = '5.31'
$ws.Range("E49").Value = '  -14.26%  '
$ws.Range("D50").Value = '0.0212'
$ws.Range("E50").Value = '  -10.64%  '
$ws.Range("D51").Value = '17.57'
$ws.Range("E51").Value = '  -14.65%  '

$priceVolRange.NumberFormat = "General"
$priceVolRange.Style = "Normal"
